# Commit: "validate button adn download works"
# Update the sample data rows (2 and 3) with the new "validated" values,
# and remove the old row 4 entirely (the data set went from 3 data rows
# down to 2 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: overwrite with the new record ---
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = 354321
$ws.Range("C2").Value = "2022-2"
$ws.Range("D2").Value = 324354
$ws.Range("E2").Value = "VALLe DORADO"
$ws.Range("F2").Value = 35435
$ws.Range("G2").Value = "ensenada"
$ws.Range("H2").Value = 54354
$ws.Range("I2").Value = "JUDITH"
$ws.Range("J2").Value = "LUNA"
$ws.Range("K2").Value = "SERRANO"
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = "Femenino"
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = "Doctorado"
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = "SAUZAL"
$ws.Range("T2").Value = "MEXICO"
$ws.Range("U2").Value = "BAJA CALIFORNIA"
$ws.Range("V2").Value = "ESPAÑOL"
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 1
$ws.Range("Y2").Value = "Estancias Sabaticas"
$ws.Range("Z2").Value = "judith"

# --- Row 3: overwrite with the new record ---
$ws.Range("A3").Value = 6
$ws.Range("B3").Value = 9999
$ws.Range("C3").Value = "periodo"
$ws.Range("D3").Value = 62155
$ws.Range("E3").Value = "destino"
$ws.Range("F3").Value = 5435435
$ws.Range("G3").Value = "unidaddd"
$ws.Range("H3").Value = 4524
$ws.Range("I3").Value = "fernanada"
$ws.Range("J3").Value = "alcala"
$ws.Range("K3").Value = "perrez"
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = "Femenino"
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = "Licenciatura"
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = "medicina"
$ws.Range("T3").Value = "mexico"
$ws.Range("U3").Value = "sonora"
$ws.Range("V3").Value = "español"
$ws.Range("W3").Value = 3
$ws.Range("X3").Value = 1
$ws.Range("Y3").Value = "Estancia de Investigacion"
$ws.Range("Z3").Value = "judith"

# --- Row 4 no longer exists in the validated data; remove it ---
$ws.Rows.Item(4).Delete()
